# The "Metadata" sheet lists CodeSystem properties as Property/Value pairs.
# This change (per the commit) refreshes the generation Date and adds a new
# "Jurisdiction" property row (currently left empty), right after "Contact".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 8 = "Date" / <generation timestamp>
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"

# Row 10 = "Contact", row 11 currently = "Description". Insert a blank row
# at position 11 so "Description" (and everything below it) shifts down by
# one, then fill the new row in with the "Jurisdiction" property.
$ws.Rows.Item(11).Insert()

# The freshly inserted row has no formatting yet - copy it from the row
# immediately above ("Contact") so it matches the rest of the table
# (border + top-aligned wrapped text), then overwrite with Copy's values.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
